$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 40/41 content swap (coins reordered: internet-computer now ranked 39, first-digital-usd ranked 40) ---
$ws.Range("A40").Value = "internet-computer"
$ws.Range("B40").Value = "Internet Computer"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "482076964.6058819000000000"
$ws.Range("E40").Value = "https://www.dfinityexplorer.org/#/"

$ws.Range("A41").Value = "first-digital-usd"
$ws.Range("B41").Value = "First Digital USD"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2588088088.1678160000000000"
$ws.Range("E41").Value = ""

# --- marketCapUsd (column F) refresh for all rows ---
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "1641207194690.0496751118614856"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "220337252974.8080123838024198"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "144586473851.1634674818780778"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "125766898323.9507969950387880"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "87094124885.8377206369660350"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "64667504451.9554172585768273"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "60234813484.9852638176274960"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "25346141527.2654701328672700"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "23691787244.0958637037073023"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "22064228294.1331864591042335"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "16968716904.4478553023259903"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "10659986060.5154514315835009"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "9643303953.7450280165809951"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "8967880690.4364066445444603"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "8872512198.8648064386231368"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "8296295563.3303270159417646"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "7992201792.6457126562201759"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "7517552122.1259767378015129"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "7412444842.7264710448046345"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "7217873319.0679865555636077"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "6527147267.3570225263640705"
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "6242045410.1634044404447100"
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Value = "6140514860.9344141927812286"
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = "6111778234.6691275847676367"
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "6032029266.7718930453094159"
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "5481899749.2262106400000000"
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "5364602397.9696044243195854"
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "5247501502.8275667023314166"
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "5237699589.2823942380839243"
$ws.Range("F31").NumberFormat = "@"
$ws.Range("F31").Value = "4237394546.0899766238032760"
$ws.Range("F32").NumberFormat = "@"
$ws.Range("F32").Value = "4005469081.4247078064784749"
$ws.Range("F33").NumberFormat = "@"
$ws.Range("F33").Value = "3752579971.8915311454245716"
$ws.Range("F34").NumberFormat = "@"
$ws.Range("F34").Value = "3247459454.5931240530693347"
$ws.Range("F35").NumberFormat = "@"
$ws.Range("F35").Value = "3170749559.1958473538077664"
$ws.Range("F36").NumberFormat = "@"
$ws.Range("F36").Value = "2973819546.1643635286852595"
$ws.Range("F37").NumberFormat = "@"
$ws.Range("F37").Value = "2905169116.5169184340000000"
$ws.Range("F38").NumberFormat = "@"
$ws.Range("F38").Value = "2797454361.2166424025848240"
$ws.Range("F39").NumberFormat = "@"
$ws.Range("F39").Value = "2667981155.8194357292965915"
$ws.Range("F40").NumberFormat = "@"
$ws.Range("F40").Value = "2584064946.9123296479551339"
$ws.Range("F41").NumberFormat = "@"
$ws.Range("F41").Value = "2583280325.3536634597306508"
$ws.Range("F42").NumberFormat = "@"
$ws.Range("F42").Value = "2536838011.8351593535535569"
$ws.Range("F43").NumberFormat = "@"
$ws.Range("F43").Value = "2519332462.8535555526326694"
$ws.Range("F44").NumberFormat = "@"
$ws.Range("F44").Value = "2502723147.5483450756625085"
$ws.Range("F45").NumberFormat = "@"
$ws.Range("F45").Value = "2132526451.9658020452506153"
$ws.Range("F46").NumberFormat = "@"
$ws.Range("F46").Value = "2043751583.7729966744552976"
$ws.Range("F47").NumberFormat = "@"
$ws.Range("F47").Value = "1974168246.4874274039062500"
$ws.Range("F48").NumberFormat = "@"
$ws.Range("F48").Value = "1960605547.5072636065451263"
$ws.Range("F49").NumberFormat = "@"
$ws.Range("F49").Value = "1960187670.1666764025242768"
$ws.Range("F50").NumberFormat = "@"
$ws.Range("F50").Value = "1917352673.0831533078970259"
$ws.Range("F51").NumberFormat = "@"
$ws.Range("F51").Value = "1911368024.4000523827148876"
$ws.Range("F52").NumberFormat = "@"
$ws.Range("F52").Value = "1829675548.5288383782389746"
$ws.Range("F53").NumberFormat = "@"
$ws.Range("F53").Value = "1794346876.4148975404160060"
$ws.Range("F54").NumberFormat = "@"
$ws.Range("F54").Value = "1687107422.5141181281472516"
$ws.Range("F55").NumberFormat = "@"
$ws.Range("F55").Value = "1672487758.0683583333113078"
$ws.Range("F56").NumberFormat = "@"
$ws.Range("F56").Value = "1552633403.3903020859105310"
$ws.Range("F57").NumberFormat = "@"
$ws.Range("F57").Value = "1547508740.9621073660913224"
$ws.Range("F58").NumberFormat = "@"
$ws.Range("F58").Value = "1469108246.9620492698586157"
$ws.Range("F59").NumberFormat = "@"
$ws.Range("F59").Value = "1444885033.7650688640000000"
$ws.Range("F60").NumberFormat = "@"
$ws.Range("F60").Value = "1367680428.0332902066306528"
$ws.Range("F61").NumberFormat = "@"
$ws.Range("F61").Value = "1358888696.5070987524859802"
$ws.Range("F62").NumberFormat = "@"
$ws.Range("F62").Value = "1311417512.5856234733158540"
$ws.Range("F63").NumberFormat = "@"
$ws.Range("F63").Value = "1226760960.9952483908518358"
$ws.Range("F64").NumberFormat = "@"
$ws.Range("F64").Value = "1202451897.8549062440000000"
$ws.Range("F65").NumberFormat = "@"
$ws.Range("F65").Value = "1083740388.2342325489542163"
$ws.Range("F66").NumberFormat = "@"
$ws.Range("F66").Value = "1065454198.6407229750000000"
$ws.Range("F67").NumberFormat = "@"
$ws.Range("F67").Value = "1049250389.0667149854776313"
$ws.Range("F68").NumberFormat = "@"
$ws.Range("F68").Value = "962507580.9806229712757570"
$ws.Range("F69").NumberFormat = "@"
$ws.Range("F69").Value = "958571455.1192458108830055"
$ws.Range("F70").NumberFormat = "@"
$ws.Range("F70").Value = "953591446.1190479056549185"
$ws.Range("F71").NumberFormat = "@"
$ws.Range("F71").Value = "924549006.9580034266694475"
$ws.Range("F73").NumberFormat = "@"
$ws.Range("F73").Value = "891548861.1944960000654884"
$ws.Range("F74").NumberFormat = "@"
$ws.Range("F74").Value = "882787928.7271338279235106"
$ws.Range("F75").NumberFormat = "@"
$ws.Range("F75").Value = "881726841.4225415106382148"
$ws.Range("F76").NumberFormat = "@"
$ws.Range("F76").Value = "876452745.4987871402788793"
$ws.Range("F77").NumberFormat = "@"
$ws.Range("F77").Value = "875455769.9508351999554378"
$ws.Range("F78").NumberFormat = "@"
$ws.Range("F78").Value = "843164617.3574918256825774"
$ws.Range("F79").NumberFormat = "@"
$ws.Range("F79").Value = "825090255.4090631000000000"
$ws.Range("F80").NumberFormat = "@"
$ws.Range("F80").Value = "815748307.4967140905742437"
$ws.Range("F81").NumberFormat = "@"
$ws.Range("F81").Value = "801762690.4410470682681000"
$ws.Range("F82").NumberFormat = "@"
$ws.Range("F82").Value = "792133992.3010214393520000"
$ws.Range("F83").NumberFormat = "@"
$ws.Range("F83").Value = "781963907.7790665377340163"
$ws.Range("F84").NumberFormat = "@"
$ws.Range("F84").Value = "767207208.6997688288921932"
$ws.Range("F85").NumberFormat = "@"
$ws.Range("F85").Value = "692716043.1744293225867631"
$ws.Range("F86").NumberFormat = "@"
$ws.Range("F86").Value = "689328171.0181693643410743"
$ws.Range("F87").NumberFormat = "@"
$ws.Range("F87").Value = "682073294.9507409179461589"
$ws.Range("F88").NumberFormat = "@"
$ws.Range("F88").Value = "679515144.1914802995852663"
$ws.Range("F89").NumberFormat = "@"
$ws.Range("F89").Value = "674628377.5776187305721227"
$ws.Range("F90").NumberFormat = "@"
$ws.Range("F90").Value = "664200063.8926383444105000"
$ws.Range("F91").NumberFormat = "@"
$ws.Range("F91").Value = "658079319.4493791258409206"
$ws.Range("F92").NumberFormat = "@"
$ws.Range("F92").Value = "639076930.0121098553281614"
$ws.Range("F93").NumberFormat = "@"
$ws.Range("F93").Value = "635134845.9482837969501684"
$ws.Range("F94").NumberFormat = "@"
$ws.Range("F94").Value = "629929521.8492605763358876"
$ws.Range("F95").NumberFormat = "@"
$ws.Range("F95").Value = "618044828.8354813040505744"
$ws.Range("F96").NumberFormat = "@"
$ws.Range("F96").Value = "604688115.9882250625751175"
$ws.Range("F97").NumberFormat = "@"
$ws.Range("F97").Value = "604192008.2020432893830351"
$ws.Range("F98").NumberFormat = "@"
$ws.Range("F98").Value = "584083082.9055908344448669"
$ws.Range("F99").NumberFormat = "@"
$ws.Range("F99").Value = "578906400.6221349658349037"
$ws.Range("F100").NumberFormat = "@"
$ws.Range("F100").Value = "569217007.7225569586693963"
$ws.Range("F101").NumberFormat = "@"
$ws.Range("F101").Value = "566159850.6601782715478274"
